# Apply edits described by the diff:
# - In sheet "ventas" (first sheet), rename header B1 from "ventas_totales" to "ingresos_totales"
# - Select cell B2 on the "ventas" sheet and make it the active/selected sheet (tab selected)
# - Column B on "ventas" widens automatically to fit new, longer header text (bestFit)

$wb = $excel.ActiveWorkbook

# Rename the header in the "ventas" worksheet
$wsVentas = $wb.Worksheets.Item("ventas")
$wsVentas.Range("B1").Value = "ingresos_totales"

# Autofit column B so the width matches the new (longer) header text
$wsVentas.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Make "ventas" the active sheet and select cell B2 on it
$wsVentas.Activate()
$wsVentas.Range("B2").Select() | Out-Null
